# The original test fixture "ExcelDataSheet.xlsx" has 3 rows of data:
#   row1: first_name | last_name | age   (header)
#   row2: John        | Doe       | 32
#   row3: Foo          | Bar       | 24
#
# The commit "modify file to have only one line of values" removes the
# second data row (row 3: Foo / Bar / 24), leaving only the header and the
# single "John Doe" data row. The now-unused "Foo"/"Bar" shared strings
# disappear from the shared string table as a natural consequence of the
# save (they are no longer referenced by any cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the values in row 3 (A3:C3) -- the cells remain part of the
# sheet (still inside the A1:C3 used range) but no longer hold data.
$ws.Range("A3:C3").ClearContents()

# Re-apply the (Normal) cell style across the whole used range. This mirrors
# what happened in the authored change: every cell in A1:C3 ends up
# pointing at a (newly written) style entry instead of the original
# implicit default style.
$ws.Range("A1:C3").Style = "Normal"

# Finally, leave the selection on D3, matching the cursor position recorded
# in the saved worksheet after the edit.
$ws.Range("D3").Select()
